$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 870.913
$ws.Range("J17").Value = 785.2727
$ws.Range("L17").Value = 2355.8181
$ws.Range("N17").Value = -2691.8181
$ws.Range("H53").Value = 848
$ws.Range("I53").Value = 284
$ws.Range("J53").Value = 1331.4286
$ws.Range("K53").Value = 284
$ws.Range("L53").Value = 1331.4286
$ws.Range("M53").Value = 353
$ws.Range("N53").Value = -2605.4286
$ws.Range("H113").Value = 7630.8945
$ws.Range("I113").Value = 7561.4
$ws.Range("J113").Value = 8127.2856
$ws.Range("K113").Value = 7561.4
$ws.Range("L113").Value = 8127.2856
$ws.Range("M113").Value = -4307.4
$ws.Range("N113").Value = -14635.2856
$ws.Range("H116").Value = 18189
$ws.Range("J116").Value = 4139.2
$ws.Range("L116").Value = 4139.2
$ws.Range("N116").Value = -11023.2
$ws.Range("H125").Value = 9148.333000000001
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 9148.333000000001
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -87254.997
$ws.Range("H132").Value = 2133.1924
$ws.Range("I132").Value = 1813.9445
$ws.Range("J132").Value = 2851.5
$ws.Range("K132").Value = 5441.833500000001
$ws.Range("L132").Value = 8554.5
$ws.Range("M132").Value = -2911.833500000001
$ws.Range("N132").Value = -13614.5
$ws.Range("H135").Value = 837.25
$ws.Range("J135").Value = 1434.5834
$ws.Range("L135").Value = 12911.2506
$ws.Range("N135").Value = -17981.2506
$ws.Range("H138").Value = 2524.03
$ws.Range("I138").Value = 1108.0769
$ws.Range("J138").Value = 2735.6091
$ws.Range("K138").Value = 3324.2307
$ws.Range("L138").Value = 8206.827300000001
$ws.Range("M138").Value = 1815.7693
$ws.Range("N138").Value = -18486.8273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14926356
$ws.Range("I32").Value = 15152391
$ws.Range("K32").Value = 15152391
$ws.Range("M32").Value = -15152104
$ws.Range("H61").Value = 4252.875
$ws.Range("I61").Value = 3950
$ws.Range("K61").Value = 3950
$ws.Range("M61").Value = -3738
$ws.Range("H92").Value = 54519.855
$ws.Range("J92").Value = 54519.855
$ws.Range("L92").Value = 54519.855
$ws.Range("N92").Value = -59511.855
$ws.Range("H97").Value = 1903.75
$ws.Range("I97").Value = 1368.6666
$ws.Range("J97").Value = 2925.2727
$ws.Range("K97").Value = 1368.6666
$ws.Range("L97").Value = 2925.2727
$ws.Range("M97").Value = -872.6666
$ws.Range("N97").Value = -3917.2727
$ws.Range("H110").Value = 1971.7084
$ws.Range("I110").Value = 1729.619
$ws.Range("K110").Value = 1729.619
$ws.Range("M110").Value = 315.3810000000001
$ws.Range("H136").Value = 4252.875
$ws.Range("I136").Value = 3950
$ws.Range("K136").Value = 11850
$ws.Range("M136").Value = -9300
$ws.Range("H139").Value = 38600
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 92496.336
$ws.Range("J116").Value = 92496.336
$ws.Range("L116").Value = 92496.336
$ws.Range("N116").Value = -101674.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2846.3655
$ws.Range("I58").Value = 2675.261
$ws.Range("K58").Value = 2675.261
$ws.Range("M58").Value = -2472.261
$ws.Range("H100").Value = 60041.125
$ws.Range("J100").Value = 60041.125
$ws.Range("L100").Value = 60041.125
$ws.Range("N100").Value = -62205.125
$ws.Range("H105").Value = 2216.65
$ws.Range("J105").Value = 3083.3333
$ws.Range("L105").Value = 3083.3333
$ws.Range("N105").Value = -6577.3333
$ws.Range("H134").Value = 1971.9375
$ws.Range("I134").Value = 1432
$ws.Range("K134").Value = 4296
$ws.Range("M134").Value = -1761
$ws.Range("H136").Value = 2846.3655
$ws.Range("I136").Value = 2675.261
$ws.Range("K136").Value = 8025.782999999999
$ws.Range("M136").Value = -5475.782999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 929.04346
$ws.Range("I7").Value = 184
$ws.Range("K7").Value = 552
$ws.Range("M7").Value = -440
$ws.Range("H11").Value = 207161.97
$ws.Range("I11").Value = 238126.67
$ws.Range("J11").Value = 125879.625
$ws.Range("K11").Value = 714380.01
$ws.Range("L11").Value = 377638.875
$ws.Range("M11").Value = -714240.01
$ws.Range("N11").Value = -377918.875
$ws.Range("H14").Value = 9529.666999999999
$ws.Range("I14").Value = 9529.666999999999
$ws.Range("K14").Value = 28589.001
$ws.Range("M14").Value = -28416.001
$ws.Range("H52").Value = 2920.8572
$ws.Range("J52").Value = 2920.8572
$ws.Range("L52").Value = 8762.571599999999
$ws.Range("N52").Value = -9294.571599999999
$ws.Range("H98").Value = 1815
$ws.Range("I98").Value = 1815
$ws.Range("K98").Value = 5445
$ws.Range("M98").Value = -3947
$ws.Range("H120").Value = 21845.924
$ws.Range("I120").Value = 11332.333
$ws.Range("K120").Value = 33996.999
$ws.Range("M120").Value = -29158.999
$ws.Range("H122").Value = 878.4
$ws.Range("J122").Value = 873
$ws.Range("L122").Value = 7857
$ws.Range("N122").Value = -12757

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 50874.5
$ws.Range("I113").Value = 34499.668
$ws.Range("K113").Value = 34499.668
$ws.Range("M113").Value = -32329.668
$ws.Range("H134").Value = 78572
$ws.Range("J134").Value = 78572
$ws.Range("L134").Value = 235716
$ws.Range("N134").Value = -240786

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6099.4346
$ws.Range("I7").Value = 6869.2
$ws.Range("K7").Value = 6869.2
$ws.Range("M7").Value = -6757.2
$ws.Range("H36").Value = 69993.5
$ws.Range("J36").Value = 69993.5
$ws.Range("L36").Value = 69993.5
$ws.Range("N36").Value = -71117.5
$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 3500
$ws.Range("K61").Value = 3500
$ws.Range("M61").Value = -3298
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("M113").Value = -1330
$ws.Range("H124").Value = 99990
$ws.Range("J124").Value = 99990
$ws.Range("L124").Value = 99990
$ws.Range("N124").Value = -109810
$ws.Range("H126").Value = 6099.4346
$ws.Range("I126").Value = 6869.2
$ws.Range("K126").Value = 20607.6
$ws.Range("M126").Value = -18137.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8758.6
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 8758.6
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -9538.6
$ws.Range("H59").Value = 49999
$ws.Range("J59").Value = 49999
$ws.Range("L59").Value = 49999
$ws.Range("N59").Value = -51475
$ws.Range("H122").Value = 10348
$ws.Range("I122").Value = 11257.333
$ws.Range("J122").Value = 8529.333000000001
$ws.Range("K122").Value = 33771.999
$ws.Range("L122").Value = 25587.999
$ws.Range("M122").Value = -31321.999
$ws.Range("N122").Value = -30487.999
$ws.Range("H126").Value = 3779.6
$ws.Range("I126").Value = 2966.3333
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 8898.999899999999
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -6428.999899999999
$ws.Range("N126").Value = -19938.5
